$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 247; existing rows 247-295 shift down to 248-296.
$ws.Rows.Item(247).Insert()

# Populate the newly inserted row 247 with the new record.
$ws.Cells.Item(247, 1).Value  = 9
$ws.Cells.Item(247, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(247, 3).Value  = "Metropolitana"
$ws.Cells.Item(247, 4).Value  = 44711
$ws.Cells.Item(247, 5).Value  = 13
$ws.Cells.Item(247, 6).Value  = 100112021
$ws.Cells.Item(247, 7).Value  = "Ají"
$ws.Cells.Item(247, 8).Value  = "Inferno"
$ws.Cells.Item(247, 9).Value  = "Primera"
$ws.Cells.Item(247, 10).Value = 79
$ws.Cells.Item(247, 11).Value = 24000
$ws.Cells.Item(247, 12).Value = 25000
$ws.Cells.Item(247, 13).Value = 24494
$ws.Cells.Item(247, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(247, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(247, 16).Value = 2041
$ws.Cells.Item(247, 17).Value = 12
$ws.Cells.Item(247, 18).Value = "Hortaliza"
